function Wrap-Xml($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: Fill in the empty "Precision and Recall" answer cell (table 2,
# row 3, column 2) with the explanatory paragraphs.
# ---------------------------------------------------------------------------
$pPrCommon = '<w:pPr><w:widowControl w:val="0"/><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:color w:val="999999"/></w:rPr></w:pPr>'

# --- paragraph 1 (existing paragraph: add justification + the first sentence)
$t = $d.Tables.Item(2)
$cell = $t.Cell(3, 2)
$p1 = $cell.Range.Paragraphs.Item(1)
$p1r = $p1.Range
$p1.Format.Alignment = 3
$p1r.InsertBefore("Precision and recall measure the model such that we can understand how the model performs for an individual class, as well as how it performs across classes.")

$cellA = $d.Tables.Item(2).Cell(3, 2)
$p1b = $cellA.Range.Paragraphs.Item(1)
$p1b.Range.Font.Color = 0x999999

# --- paragraph bodies 2 .. 8 (each appended via InsertParagraphAfter + InsertXML)
$bodies = @(
  '',
  ('<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="999999"/></w:rPr><w:t>Model precision</w:t></w:r>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t xml:space="preserve"> measures the percentage of correct predictions against total number of predictions.</w:t></w:r>'),
  '',
  ('<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="999999"/></w:rPr><w:t xml:space="preserve">Model </w:t></w:r>' + `
   '<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="999999"/></w:rPr><w:t>recall</w:t></w:r>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t xml:space="preserve"> measures the percentage of correct</w:t></w:r>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t>ly</w:t></w:r>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t>identified instances</w:t></w:r>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t xml:space="preserve"> against total </w:t></w:r>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t>possible instances</w:t></w:r>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t>.</w:t></w:r>'),
  '',
  ('<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="999999"/></w:rPr><w:t>Precision</w:t></w:r>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t xml:space="preserve"> = TP/(TP+FP) = 10</w:t></w:r>' + `
   '<w:proofErr w:type="gramStart"/>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t>/(</w:t></w:r>' + `
   '<w:proofErr w:type="gramEnd"/>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t>10+2) = 10/12 ~= 0,83</w:t></w:r>'),
  ('<w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="999999"/></w:rPr><w:t>Recall</w:t></w:r>' + `
   '<w:r><w:rPr><w:color w:val="999999"/></w:rPr><w:t xml:space="preserve"> = TP/(TP+FN) = 10(10+0) = 10/10 = 1</w:t></w:r>')
)

foreach ($body in $bodies) {
    $cellX = $d.Tables.Item(2).Cell(3, 2)
    $lastIdx = $cellX.Range.Paragraphs.Count
    $lastPara = $cellX.Range.Paragraphs.Item($lastIdx)
    $newPara = $lastPara.Range.InsertParagraphAfter()

    $cellY = $d.Tables.Item(2).Cell(3, 2)
    $newIdx = $cellY.Range.Paragraphs.Count
    $newRange = $cellY.Range.Paragraphs.Item($newIdx).Range
    $newRange.InsertXML((Wrap-Xml ('<w:p>' + $pPrCommon + $body + '</w:p>')))
}
